# Generate Report for Handoff
# Adds two new files (6e6f5a0a-... and 6edcb7bb-...) to the Overview,
# zh-cn and de-de tables, mirroring the existing "In Translation" rows
# but in the "Ready for handoff" state.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d7282bd5818541c6c959d3095f75540b8cd8a0b/e2e/"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$r = 4
$ws.Range("A$r").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$ws.Hyperlinks.Add($ws.Range("B$r"), ($ghBase + "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"), [Type]::Missing, [Type]::Missing, "e2e\6e6f5a0a-f847-4138-9cfe-7cdb61058920.md") | Out-Null
$ws.Range("B$r").Style = "HyperLink"
$ws.Range("C$r").Value = ".md"
$ws.Range("E$r").Value = "Ready for handoff"
$ws.Range("F$r").Value = "Ready for handoff"
$ws.Range("G$r").Value = "2016-12-16 08:11:22"
$ws.Range("G$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$lo.ListRows.Add() | Out-Null
$r = 5
$ws.Range("A$r").Value = "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md"
$ws.Hyperlinks.Add($ws.Range("B$r"), ($ghBase + "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md"), [Type]::Missing, [Type]::Missing, "e2e\6edcb7bb-6a62-49b9-b5d2-6afc280672df.md") | Out-Null
$ws.Range("B$r").Style = "HyperLink"
$ws.Range("C$r").Value = ".md"
$ws.Range("E$r").Value = "Ready for handoff"
$ws.Range("F$r").Value = "Ready for handoff"
$ws.Range("G$r").Value = "2016-12-16 08:11:22"
$ws.Range("G$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$r = 4
$ws.Hyperlinks.Add($ws.Range("A$r"), ($ghBase + "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"), [Type]::Missing, [Type]::Missing, "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md") | Out-Null
$ws.Range("A$r").Style = "HyperLink"
$ws.Range("B$r").Value = ".md"
$ws.Range("C$r").Value = "Ready for handoff"
$ws.Range("D$r").Value = "e2e"
$ws.Range("E$r").Value = "ht"
$ws.Range("F$r").Value = "False"
$ws.Range("G$r").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.zh-cn.xlf"
$ws.Range("H$r").Value = "2016-12-16 08:11:09"
$ws.Range("H$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L$r").Value = "0001-01-01 00:00:00"
$ws.Range("L$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("O$r").Value = "True"
$ws.Range("Q$r").Value = "False"

$lo.ListRows.Add() | Out-Null
$r = 5
$ws.Hyperlinks.Add($ws.Range("A$r"), ($ghBase + "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md"), [Type]::Missing, [Type]::Missing, "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md") | Out-Null
$ws.Range("A$r").Style = "HyperLink"
$ws.Range("B$r").Value = ".md"
$ws.Range("C$r").Value = "Ready for handoff"
$ws.Range("D$r").Value = "e2e"
$ws.Range("E$r").Value = "ht"
$ws.Range("F$r").Value = "False"
$ws.Range("G$r").Value = "6edcb7bb-6a62-49b9-b5d2-6afc280672df.7d2fde898036ef06dfabdcf653b03f9c64799115.zh-cn.xlf"
$ws.Range("H$r").Value = "2016-12-16 08:11:09"
$ws.Range("H$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L$r").Value = "0001-01-01 00:00:00"
$ws.Range("L$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("O$r").Value = "True"
$ws.Range("Q$r").Value = "False"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$r = 4
$ws.Hyperlinks.Add($ws.Range("A$r"), ($ghBase + "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"), [Type]::Missing, [Type]::Missing, "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md") | Out-Null
$ws.Range("A$r").Style = "HyperLink"
$ws.Range("B$r").Value = ".md"
$ws.Range("C$r").Value = "Ready for handoff"
$ws.Range("D$r").Value = "e2e"
$ws.Range("E$r").Value = "ht"
$ws.Range("F$r").Value = "False"
$ws.Range("G$r").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.de-de.xlf"
$ws.Range("H$r").Value = "2016-12-16 08:11:22"
$ws.Range("H$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L$r").Value = "0001-01-01 00:00:00"
$ws.Range("L$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("O$r").Value = "True"
$ws.Range("Q$r").Value = "False"

$lo.ListRows.Add() | Out-Null
$r = 5
$ws.Hyperlinks.Add($ws.Range("A$r"), ($ghBase + "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md"), [Type]::Missing, [Type]::Missing, "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md") | Out-Null
$ws.Range("A$r").Style = "HyperLink"
$ws.Range("B$r").Value = ".md"
$ws.Range("C$r").Value = "Ready for handoff"
$ws.Range("D$r").Value = "e2e"
$ws.Range("E$r").Value = "ht"
$ws.Range("F$r").Value = "False"
$ws.Range("G$r").Value = "6edcb7bb-6a62-49b9-b5d2-6afc280672df.7d2fde898036ef06dfabdcf653b03f9c64799115.de-de.xlf"
$ws.Range("H$r").Value = "2016-12-16 08:11:22"
$ws.Range("H$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L$r").Value = "0001-01-01 00:00:00"
$ws.Range("L$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("O$r").Value = "True"
$ws.Range("Q$r").Value = "False"

Write-Output "done"
